# Generate Report for Handback
# Updates the handoff/handback timestamps and priority value that are
# refreshed each time the handback status report is regenerated.

$wb = $excel.ActiveWorkbook

# --- "Overview" sheet: Latest HO Xliff Generate Date (column G) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-13 20:21:08"
$wsOverview.Range("G4").Value = "2016-08-13 20:21:08"

# --- "zh-cn" sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
# Priority column (E): ht -> mt
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E4").Value = "mt"
# Correspond Handoff Datetime (H)
$wsZhCn.Range("H2").Value = "2016-08-13 20:20:56"
$wsZhCn.Range("H4").Value = "2016-08-13 20:20:56"
# Correspond Handback DateTime (K)
$wsZhCn.Range("K2").Value = "2016-08-13 20:21:27"
$wsZhCn.Range("K4").Value = "2016-08-13 20:21:27"

# --- "de-de" sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
# Priority column (E): ht -> mt
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E4").Value = "mt"
# Correspond Handoff Datetime (H)
$wsDeDe.Range("H2").Value = "2016-08-13 20:21:08"
$wsDeDe.Range("H4").Value = "2016-08-13 20:21:08"
# Correspond Handback DateTime (K)
$wsDeDe.Range("K2").Value = "2016-08-13 20:21:37"
$wsDeDe.Range("K4").Value = "2016-08-13 20:21:37"
